$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E..I hold five consecutive "twelve months ended" periods.
# The oldest period (column E) is dropped, every period shifts one column
# to the left, and the new period"s data lands in column I.
$dataCols = @(5, 6, 7, 8, 9)   # E, F, G, H, I (numeric column indices)

$headerRows = @(8, 25)
$newHeader = "دوازده ماهه منتهی به 1401/12"
foreach ($r in $headerRows) {
    for ($i = 0; $i -lt 4; $i++) {
        $srcCol = $dataCols[$i + 1]
        $dstCol = $dataCols[$i]
        $ws.Cells.Item($r, $dstCol).Value2 = $ws.Cells.Item($r, $srcCol).Value2
    }
    $ws.Cells.Item($r, $dataCols[4]).Value2 = $newHeader
}

$newPeriodValues = @{
    10 = 0
    11 = 0
    12 = 0
    13 = 86979
    14 = 0
    15 = 0
    16 = 0
    17 = 37030
    18 = 107678
    19 = 0
    20 = 110019
    21 = 341706
    27 = 60
    28 = 228
    29 = "-"
    30 = "-"
}

foreach ($r in $newPeriodValues.Keys) {
    for ($i = 0; $i -lt 4; $i++) {
        $srcCol = $dataCols[$i + 1]
        $dstCol = $dataCols[$i]
        $ws.Cells.Item($r, $dstCol).Value2 = $ws.Cells.Item($r, $srcCol).Value2
    }
    $ws.Cells.Item($r, $dataCols[4]).Value2 = $newPeriodValues[$r]
}
